# Generate Report for Handoff
# Update status of b5888c65-814d-43b5-a0df-dbc30d19cc59.md to "Ready for handoff"
# and set new handoff datetimes on the zh-cn and de-de sheets / Overview.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Overview sheet: row 3 corresponds to b5888c65-814d-43b5-a0df-dbc30d19cc59.md
$overview.Range("B3").Value = "Ready for handoff"
$overview.Range("C3").Value = "Ready for handoff"

# zh-cn sheet: row 3 -> Status + Latest Handoff Datetime
$zhcn.Range("B3").Value = "Ready for handoff"
$zhcn.Range("D3").Value = "2016-03-01 05:32:59"

# de-de sheet: row 3 -> Status + Latest Handoff Datetime
$dede.Range("B3").Value = "Ready for handoff"
$dede.Range("D3").Value = "2016-03-01 05:33:10"
